$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the referenced cell addresses from row 24 to row 25 (new init for scenario 29)
$ws.Range("D5").Value  = "A25"
$ws.Range("D6").Value  = "B25"
$ws.Range("D7").Value  = "C25"
$ws.Range("D8").Value  = "G25"
$ws.Range("D9").Value  = "H25"
$ws.Range("D10").Value = "I25"
$ws.Range("D11").Value = "J25"

# Move the active selection to D12
$ws.Range("D12").Select()
